# Auto-generated Excel COM-interop script
# Applies the cell-level numeric updates described in the commit diff
# for the 'Seraph_Profits' style leveling tables across 8 sheets (ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 8124.5
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H86").Value = 5145
$ws.Range("I86").Value = 5899
$ws.Range("J86").Value = 4977.4443
$ws.Range("K86").Value = 5899
$ws.Range("L86").Value = 4977.4443
$ws.Range("M86").Value = -4776
$ws.Range("N86").Value = -7223.4443

$ws.Range("H89").Value = 5145
$ws.Range("I89").Value = 5899
$ws.Range("J89").Value = 4977.4443
$ws.Range("K89").Value = 29495
$ws.Range("L89").Value = 24887.2215
$ws.Range("M89").Value = -23879
$ws.Range("N89").Value = -36119.2215

$ws.Range("H113").Value = 2247.5
$ws.Range("I113").Value = 2330
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2330
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 924
$ws.Range("N113").Value = -8508

$ws.Range("H116").Value = 4000.8333
$ws.Range("I116").Value = 3801
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3801
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -359
$ws.Range("N116").Value = -11884

$ws.Range("H132").Value = 1956.9231
$ws.Range("I132").Value = 2140.4546
$ws.Range("J132").Value = 947.5
$ws.Range("K132").Value = 6421.3638
$ws.Range("L132").Value = 2842.5
$ws.Range("M132").Value = -3891.3638
$ws.Range("N132").Value = -7902.5

$ws.Range("H141").Value = 8273.75
$ws.Range("J141").Value = 8500
$ws.Range("L141").Value = 25500
$ws.Range("N141").Value = -35860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1380
$ws.Range("I21").Value = 1006.6667
$ws.Range("J21").Value = 2500
$ws.Range("K21").Value = 1006.6667
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = -632.6667
$ws.Range("N21").Value = -3248

$ws.Range("H30").Value = 15199.5
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H45").Value = 1637
$ws.Range("I45").Value = 1637
$ws.Range("K45").Value = 1637
$ws.Range("M45").Value = -1260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 341.2
$ws.Range("I37").Value = 301.75
$ws.Range("J37").Value = 499
$ws.Range("K37").Value = 301.75
$ws.Range("L37").Value = 499
$ws.Range("M37").Value = -164.75
$ws.Range("N37").Value = -773

$ws.Range("H105").Value = 4762.7896
$ws.Range("I105").Value = 3486.2
$ws.Range("J105").Value = 9550
$ws.Range("K105").Value = 3486.2
$ws.Range("L105").Value = 9550
$ws.Range("M105").Value = -1739.2
$ws.Range("N105").Value = -13044

$ws.Range("H134").Value = 2026.138
$ws.Range("I134").Value = 1727.5
$ws.Range("J134").Value = 3459.6
$ws.Range("K134").Value = 5182.5
$ws.Range("L134").Value = 10378.8
$ws.Range("M134").Value = -2647.5
$ws.Range("N134").Value = -15448.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3678
$ws.Range("I132").Value = 3678
$ws.Range("K132").Value = 11034
$ws.Range("M132").Value = -8504

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 879.5
$ws.Range("J12").Value = 885.1429000000001
$ws.Range("L12").Value = 2655.4287
$ws.Range("N12").Value = -3001.4287

$ws.Range("H97").Value = 149.75
$ws.Range("I97").Value = 149.66667
$ws.Range("J97").Value = 150
$ws.Range("K97").Value = 449.00001
$ws.Range("L97").Value = 450
$ws.Range("M97").Value = 46.99998999999997
$ws.Range("N97").Value = -1442

$ws.Range("H129").Value = 2279.7
$ws.Range("I129").Value = 752.8
$ws.Range("J129").Value = 3806.6
$ws.Range("K129").Value = 2258.4
$ws.Range("L129").Value = 11419.8
$ws.Range("M129").Value = 2741.6
$ws.Range("N129").Value = -21419.8

$ws.Range("H139").Value = 2179.8
$ws.Range("I139").Value = 1600
$ws.Range("K139").Value = 4800
$ws.Range("M139").Value = 340

$ws.Range("H140").Value = 350234.5
$ws.Range("I140").Value = 350234.5
$ws.Range("K140").Value = 1050703.5
$ws.Range("M140").Value = -1045523.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 263.3158
$ws.Range("I2").Value = 64.5
$ws.Range("J2").Value = 484.22223
$ws.Range("K2").Value = 64.5
$ws.Range("L2").Value = 484.22223
$ws.Range("M2").Value = 48.5
$ws.Range("N2").Value = -710.2222300000001

$ws.Range("H43").Value = 2250
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2250
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2250
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2552

$ws.Range("H57").Value = 29965.666
$ws.Range("J57").Value = 29965.666
$ws.Range("L57").Value = 29965.666
$ws.Range("N57").Value = -31605.666

$ws.Range("H80").Value = 150
$ws.Range("I80").Value = 150
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 150
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 848
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 150
$ws.Range("I83").Value = 150
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 750
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 4242
$ws.Range("N83").ClearContents()

$ws.Range("H122").Value = 65444.438
$ws.Range("I122").Value = 2731.6924
$ws.Range("J122").Value = 337199.66
$ws.Range("K122").Value = 8195.0772
$ws.Range("L122").Value = 1011598.98
$ws.Range("M122").Value = -5745.0772
$ws.Range("N122").Value = -1016498.98

$ws.Range("H132").Value = 4067.0715
$ws.Range("I132").Value = 2456.111
$ws.Range("K132").Value = 7368.333
$ws.Range("M132").Value = -4838.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2398
$ws.Range("I7").Value = 2535.125
$ws.Range("J7").Value = 2123.75
$ws.Range("K7").Value = 2535.125
$ws.Range("L7").Value = 2123.75
$ws.Range("M7").Value = -2423.125
$ws.Range("N7").Value = -2347.75

$ws.Range("H40").Value = 3065.0908
$ws.Range("I40").Value = 2559.5715
$ws.Range("K40").Value = 2559.5715
$ws.Range("M40").Value = -2423.5715

$ws.Range("H55").Value = 233.80952
$ws.Range("I55").Value = 243.5
$ws.Range("K55").Value = 243.5
$ws.Range("M55").Value = -70.5

$ws.Range("H122").Value = 7859.6665
$ws.Range("I122").Value = 8877.223
$ws.Range("J122").Value = 6333.3335
$ws.Range("K122").Value = 26631.669
$ws.Range("L122").Value = 19000.0005
$ws.Range("M122").Value = -24181.669
$ws.Range("N122").Value = -23900.0005

$ws.Range("H126").Value = 2398
$ws.Range("I126").Value = 2535.125
$ws.Range("J126").Value = 2123.75
$ws.Range("K126").Value = 7605.375
$ws.Range("L126").Value = 6371.25
$ws.Range("M126").Value = -5135.375
$ws.Range("N126").Value = -11311.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 21000
$ws.Range("J104").Value = 21000
$ws.Range("L104").Value = 21000
$ws.Range("N104").Value = -27988

$ws.Range("H132").Value = 1602.25
$ws.Range("I132").Value = 1468
$ws.Range("K132").Value = 4404
$ws.Range("M132").Value = -1874
